$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column at I (old I..M shift right by one, to J..N)
$ws.Range("I1").EntireColumn.Insert()
# Remove the now-duplicated blank column that resulted from the shift (keeps N..R aligned)
$ws.Range("O1").EntireColumn.Delete()
